$wb = $excel.ActiveWorkbook

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2287.7144
$ws.Range("I62").Value = 2452.3333
$ws.Range("K62").Value = 2452.3333
$ws.Range("M62").Value = -1828.3333

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 2287.7144
$ws.Range("I65").Value = 2452.3333
$ws.Range("K65").Value = 12261.6665
$ws.Range("M65").Value = -9141.666499999999

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3941.0588
$ws.Range("I76").Value = 3833.3333
$ws.Range("J76").Value = 4062.25
$ws.Range("K76").Value = 3833.3333
$ws.Range("L76").Value = 4062.25
$ws.Range("M76").Value = -3518.3333
$ws.Range("N76").Value = -4692.25

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3941.0588
$ws.Range("I79").Value = 3833.3333
$ws.Range("J79").Value = 4062.25
$ws.Range("K79").Value = 3833.3333
$ws.Range("L79").Value = 4062.25
$ws.Range("M79").Value = -2741.3333
$ws.Range("N79").Value = -6246.25

# ALC row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4877.1816
$ws.Range("I86").Value = 4417.25
$ws.Range("J86").Value = 5140
$ws.Range("K86").Value = 4417.25
$ws.Range("L86").Value = 5140
$ws.Range("M86").Value = -3294.25
$ws.Range("N86").Value = -7386

# ALC row 88
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 1857.2222
$ws.Range("I88").Value = 3191.4285
$ws.Range("J88").Value = 1008.1818
$ws.Range("K88").Value = 3191.4285
$ws.Range("L88").Value = 1008.1818
$ws.Range("M88").Value = -2785.4285
$ws.Range("N88").Value = -1820.1818

# ALC row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 4877.1816
$ws.Range("I89").Value = 4417.25
$ws.Range("J89").Value = 5140
$ws.Range("K89").Value = 22086.25
$ws.Range("L89").Value = 25700
$ws.Range("M89").Value = -16470.25
$ws.Range("N89").Value = -36932

# ALC row 91
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 1857.2222
$ws.Range("I91").Value = 3191.4285
$ws.Range("J91").Value = 1008.1818
$ws.Range("K91").Value = 3191.4285
$ws.Range("L91").Value = 1008.1818
$ws.Range("M91").Value = -1787.4285
$ws.Range("N91").Value = -3816.1818

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1788.1818
$ws.Range("J112").Value = 2013.3334
$ws.Range("L112").Value = 6040.0002
$ws.Range("N112").Value = -8256.0002

# ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1781.3103
$ws.Range("I135").Value = 1734.04
$ws.Range("J135").Value = 2076.75
$ws.Range("K135").Value = 15606.36
$ws.Range("L135").Value = 18690.75
$ws.Range("M135").Value = -13071.36
$ws.Range("N135").Value = -23760.75

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1860.44
$ws.Range("I138").Value = 867.4103
$ws.Range("J138").Value = 2495.328
$ws.Range("K138").Value = 2602.2309
$ws.Range("L138").Value = 7485.984
$ws.Range("M138").Value = 2537.7691
$ws.Range("N138").Value = -17765.984

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3183.71
$ws.Range("I32").Value = 3146.1716
$ws.Range("J32").Value = 6900
$ws.Range("K32").Value = 3146.1716
$ws.Range("L32").Value = 6900
$ws.Range("M32").Value = -2859.1716
$ws.Range("N32").Value = -7474

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1977997.2
$ws.Range("I45").Value = 2526979
$ws.Range("K45").Value = 2526979
$ws.Range("M45").Value = -2526602

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3752.682
$ws.Range("I61").Value = 3752.682
$ws.Range("K61").Value = 3752.682
$ws.Range("M61").Value = -3540.682

# ARM row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2000
$ws.Range("I88").Value = 2000
$ws.Range("J88").Value = 2000
$ws.Range("K88").Value = 2000
$ws.Range("L88").Value = 2000
$ws.Range("M88").Value = -1594
$ws.Range("N88").Value = -2812

# ARM row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 2000
$ws.Range("I91").Value = 2000
$ws.Range("J91").Value = 2000
$ws.Range("K91").Value = 2000
$ws.Range("L91").Value = 2000
$ws.Range("M91").Value = -596
$ws.Range("N91").Value = -4808

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3752.682
$ws.Range("I136").Value = 3752.682
$ws.Range("K136").Value = 11258.046
$ws.Range("M136").Value = -8708.045999999998

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1177.875
$ws.Range("I107").Value = 1038.8485
$ws.Range("J107").Value = 1833.2858
$ws.Range("K107").Value = 1038.8485
$ws.Range("L107").Value = 1833.2858
$ws.Range("M107").Value = 881.1514999999999
$ws.Range("N107").Value = -5673.2858

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3664.3333
$ws.Range("I134").Value = 2439.4482
$ws.Range("J134").Value = 5278.9546
$ws.Range("K134").Value = 7318.344599999999
$ws.Range("L134").Value = 15836.8638
$ws.Range("M134").Value = -4783.344599999999
$ws.Range("N134").Value = -20906.8638

# CRP row 6
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 6200120
$ws.Range("I6").Value = 5777911
$ws.Range("J6").Value = 10000000
$ws.Range("K6").Value = 5777911
$ws.Range("L6").Value = 10000000
$ws.Range("M6").Value = -5777798
$ws.Range("N6").Value = -10000226

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2237.5466
$ws.Range("I31").Value = 2053.7222
$ws.Range("J31").Value = 2407.2307
$ws.Range("K31").Value = 2053.7222
$ws.Range("L31").Value = 2407.2307
$ws.Range("M31").Value = -1758.7222
$ws.Range("N31").Value = -2997.2307

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2237.5466
$ws.Range("I34").Value = 2053.7222
$ws.Range("J34").Value = 2407.2307
$ws.Range("K34").Value = 2053.7222
$ws.Range("L34").Value = 2407.2307
$ws.Range("M34").Value = -1851.7222
$ws.Range("N34").Value = -2811.2307

# CUL row 17
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 100
$ws.Range("J17").Value = 100
$ws.Range("L17").Value = 300
$ws.Range("N17").Value = -638

# CUL row 38
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 228.54167
$ws.Range("I38").Value = 190
$ws.Range("J38").Value = 256.07144
$ws.Range("K38").Value = 570
$ws.Range("L38").Value = 768.21432
$ws.Range("M38").Value = -223
$ws.Range("N38").Value = -1462.21432

# CUL row 59
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1212.24
$ws.Range("I122").Value = 1100.1765
$ws.Range("J122").Value = 1450.375
$ws.Range("K122").Value = 9901.5885
$ws.Range("L122").Value = 13053.375
$ws.Range("M122").Value = -7451.5885
$ws.Range("N122").Value = -17953.375

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1094.9706
$ws.Range("I131").Value = 659.8333
$ws.Range("J131").Value = 1188.2142
$ws.Range("K131").Value = 1979.4999
$ws.Range("L131").Value = 3564.6426
$ws.Range("M131").Value = 3060.5001
$ws.Range("N131").Value = -13644.6426

# CUL row 134
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 4352.4136
$ws.Range("I134").Value = 2296.1904
$ws.Range("J134").Value = 9750
$ws.Range("K134").Value = 6888.5712
$ws.Range("L134").Value = 29250
$ws.Range("M134").Value = -1818.5712
$ws.Range("N134").Value = -39390

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 201012.8
$ws.Range("I113").Value = 334117
$ws.Range("J113").Value = 1356.5
$ws.Range("K113").Value = 334117
$ws.Range("L113").Value = 1356.5
$ws.Range("M113").Value = -331947
$ws.Range("N113").Value = -5696.5

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5559.8
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 5559.8
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 16679.4
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -21619.4

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4642.75
$ws.Range("I40").Value = 6874.25
$ws.Range("J40").Value = 2411.25
$ws.Range("K40").Value = 6874.25
$ws.Range("L40").Value = 2411.25
$ws.Range("M40").Value = -6738.25
$ws.Range("N40").Value = -2683.25

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1708.6364
$ws.Range("I68").Value = 2223
$ws.Range("J68").Value = 1414.7142
$ws.Range("K68").Value = 2223
$ws.Range("L68").Value = 1414.7142
$ws.Range("M68").Value = -1474
$ws.Range("N68").Value = -2912.7142

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1708.6364
$ws.Range("I71").Value = 2223
$ws.Range("J71").Value = 1414.7142
$ws.Range("K71").Value = 11115
$ws.Range("L71").Value = 7073.571
$ws.Range("M71").Value = -7371
$ws.Range("N71").Value = -14561.571

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1248.4762
$ws.Range("I82").Value = 1132.5264
$ws.Range("J82").Value = 2350
$ws.Range("K82").Value = 1132.5264
$ws.Range("L82").Value = 2350
$ws.Range("M82").Value = -771.5264
$ws.Range("N82").Value = -3072

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1248.4762
$ws.Range("I85").Value = 1132.5264
$ws.Range("J85").Value = 2350
$ws.Range("K85").Value = 1132.5264
$ws.Range("L85").Value = 2350
$ws.Range("M85").Value = 115.4736
$ws.Range("N85").Value = -4846

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 10671.714
$ws.Range("I122").Value = 12220.8
$ws.Range("J122").Value = 9811.111000000001
$ws.Range("K122").Value = 36662.39999999999
$ws.Range("L122").Value = 29433.333
$ws.Range("M122").Value = -34212.39999999999
$ws.Range("N122").Value = -34333.333

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 12829523
$ws.Range("I132").Value = 6899.55
$ws.Range("J132").Value = 26327022
$ws.Range("K132").Value = 20698.65
$ws.Range("L132").Value = 78981066
$ws.Range("M132").Value = -18168.65
$ws.Range("N132").Value = -78986126

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2596.6296
$ws.Range("I132").Value = 2385.762
$ws.Range("J132").Value = 3334.6667
$ws.Range("K132").Value = 7157.286
$ws.Range("L132").Value = 10004.0001
$ws.Range("M132").Value = -4627.286
$ws.Range("N132").Value = -15064.0001

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2085.8667
$ws.Range("I136").Value = 1331.1052
$ws.Range("J136").Value = 3389.5454
$ws.Range("K136").Value = 3993.3156
$ws.Range("L136").Value = 10168.6362
$ws.Range("M136").Value = -1443.3156
$ws.Range("N136").Value = -15268.6362
